$wb = $excel.ActiveWorkbook

# Hyperlink target URLs and handback timestamps per language sheet.
$targets = @{
    "zh-cn" = @{
        "mdBase"   = "https://github.com/OpenLocalizationTest/oltest/blob/4d5b873cced866eeb54acbaf5131e0a109a33fa7/e2e"
        "b2891065" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b22fc7655f1598836de19167e43096179a19d832/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b2891065-cfe7-4731-94c2-4c754d4bef1c.e243ccac099635fd236a26a83151e14520daf803.zh-cn.xlf"
        "c9e26659" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b22fc7655f1598836de19167e43096179a19d832/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c9e26659-4cc2-4290-ae4b-5198713a52b8.27a04f146edc21599d6f9ab6bfa524c441f7c242.zh-cn.xlf"
        "config"   = "https://github.com/OpenLocalizationTest/oltest/blob/4d5b873cced866eeb54acbaf5131e0a109a33fa7/.localization-config"
        "handback" = "2016-02-22 18:05:22"
    }
    "de-de" = @{
        "mdBase"   = "https://github.com/OpenLocalizationTest/oltest/blob/4d5b873cced866eeb54acbaf5131e0a109a33fa7/e2e"
        "b2891065" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea656b96bd9f1743c96cae91936aaba96094a729/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b2891065-cfe7-4731-94c2-4c754d4bef1c.e243ccac099635fd236a26a83151e14520daf803.de-de.xlf"
        "c9e26659" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea656b96bd9f1743c96cae91936aaba96094a729/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c9e26659-4cc2-4290-ae4b-5198713a52b8.27a04f146edc21599d6f9ab6bfa524c441f7c242.de-de.xlf"
        "config"   = "https://github.com/OpenLocalizationTest/oltest/blob/4d5b873cced866eeb54acbaf5131e0a109a33fa7/.localization-config"
        "handback" = "2016-02-22 18:05:42"
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $t = $targets[$sheetName]

    $mdFileB = "b2891065-cfe7-4731-94c2-4c754d4bef1c.md"
    $mdFileC = "c9e26659-4cc2-4290-ae4b-5198713a52b8.md"
    $xlfFileB = "b2891065-cfe7-4731-94c2-4c754d4bef1c.e243ccac099635fd236a26a83151e14520daf803." + $sheetName + ".xlf"
    $xlfFileC = "c9e26659-4cc2-4290-ae4b-5198713a52b8.27a04f146edc21599d6f9ab6bfa524c441f7c242." + $sheetName + ".xlf"

    # Status column -> the report now reflects a handback, not a pending handoff.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # "Latest Handback DateTime" (column G) now carries a real timestamp
    # instead of the zero-date placeholder.
    $ws.Range("G2").Value = $t["handback"]
    $ws.Range("G3").Value = $t["handback"]

    # Populate the new "Latest Target File" (E) / "Latest Handback File" (F)
    # columns - same files as the source / handoff columns, now mirrored
    # back because the report is generated for the handback.
    $ws.Range("E2").Value = $mdFileB
    $ws.Range("F2").Value = $xlfFileB
    $ws.Range("E3").Value = $mdFileC
    $ws.Range("F3").Value = $xlfFileC

    # Rebuild all hyperlinks in the final left-to-right / top-to-bottom order
    # so relationship ids come out as rId2..rId10 in that same order.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), ($t["mdBase"] + "/" + $mdFileB), "", "", $mdFileB) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $t["b2891065"], "", "", $xlfFileB) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), ($t["mdBase"] + "/" + $mdFileB), "", "", $mdFileB) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $t["b2891065"], "", "", $xlfFileB) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), ($t["mdBase"] + "/" + $mdFileC), "", "", $mdFileC) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $t["c9e26659"], "", "", $xlfFileC) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E3"), ($t["mdBase"] + "/" + $mdFileC), "", "", $mdFileC) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $t["c9e26659"], "", "", $xlfFileC) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $t["config"], "", "", ".localization-config") | Out-Null

    # Hyperlinks.Add resets formatting to a generic default; restore the
    # underlined "hyperlink blue" look the rest of the link cells use.
    foreach ($addr in @("A2", "C2", "E2", "F2", "A3", "C3", "E3", "F3", "A4")) {
        $ws.Range($addr).Font.Underline = $true
        $ws.Range($addr).Font.Color = 0xED9564
    }
}
